# single_session_batch.xlsx batch edit
# - Adds "Stock Deregistration File" / "Stock VMT File" rows to the IO Settings
#   block (pushing Verbose Output / Slice Tech Combo Tables down by two rows).
# - Replaces the old "Stock Deregistration" / "Stock VMT" rows (and the blank
#   row after them) in the Postproc Settings block with seven new
#   "Context ... File" rows.
# - Updates sheet view (zoom / frozen pane / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new rows right after "ZEV Requirement File" (row 24), so the
#    "Verbose Output" / "Slice Tech Combo Tables" rows move from 25-26 to
#    27-28.
# ---------------------------------------------------------------------------
$ws.Rows.Item(25).Insert()
$ws.Rows.Item(26).Insert()

# Copy formatting from row 24 (same visual block, matching borders/shading)
# into the two freshly inserted blank rows.
$ws.Rows.Item(24).Copy()
$ws.Rows.Item(25).PasteSpecial(-4122)
$ws.Rows.Item(24).Copy()
$ws.Rows.Item(26).PasteSpecial(-4122)

$ws.Range("A25").Value = "Stock Deregistration File"
$ws.Range("B25").Value = "String"
$ws.Range("C25").Value = "input_samples/reregistration_fixed_by_age.csv"

$ws.Range("A26").Value = "Stock VMT File"
$ws.Range("B26").Value = "String"
$ws.Range("C26").Value = "input_samples/annual_vmt_fixed_by_age.csv"

# ---------------------------------------------------------------------------
# 2. In the Postproc Settings block, the old "Stock Deregistration" (row 42)
#    and "Stock VMT" (row 43) rows (now shifted to 44 & 45 after the insert
#    above) plus the trailing blank row (now 46) are replaced by seven new
#    "Context ... File" rows (44-50). Insert four extra rows so the existing
#    three (44-46) plus four new ones (47-50) give the seven needed.
# ---------------------------------------------------------------------------
$ws.Rows.Item(47).Insert()
$ws.Rows.Item(47).Insert()
$ws.Rows.Item(47).Insert()
$ws.Rows.Item(47).Insert()

# Copy formatting from row 26 (style: label bold box, text value box) into
# rows 44-49; row 50 keeps the plain style already on the (formerly blank)
# row 44 donor, so style it from row 44's original blank format instead.
$ws.Rows.Item(26).Copy()
$ws.Rows.Item(44).PasteSpecial(-4122)
$ws.Rows.Item(26).Copy()
$ws.Rows.Item(45).PasteSpecial(-4122)
$ws.Rows.Item(26).Copy()
$ws.Rows.Item(46).PasteSpecial(-4122)
$ws.Rows.Item(26).Copy()
$ws.Rows.Item(47).PasteSpecial(-4122)
$ws.Rows.Item(26).Copy()
$ws.Rows.Item(48).PasteSpecial(-4122)
$ws.Rows.Item(26).Copy()
$ws.Rows.Item(49).PasteSpecial(-4122)

$ws.Range("A44").Value = "Context Criteria Cost Factors File"
$ws.Range("B44").Value = "String"
$ws.Range("C44").Value = "input_samples/context_cost_factors-criteria.csv"

$ws.Range("A45").Value = "Context SCC Cost Factors File"
$ws.Range("B45").Value = "String"
$ws.Range("C45").Value = "input_samples/context_cost_factors-scc.csv"

$ws.Range("A46").Value = "Context Powersector Emission Factors File"
$ws.Range("B46").Value = "String"
$ws.Range("C46").Value = "input_samples/context_emission_factors-powersector.csv"

$ws.Range("A47").Value = "Context Refinery Emission Factors File"
$ws.Range("B47").Value = "String"
$ws.Range("C47").Value = "input_samples/context_emission_factors-refinery.csv"

$ws.Range("A48").Value = "Context Vehicle Emission Factors File"
$ws.Range("B48").Value = "String"
$ws.Range("C48").Value = "input_samples/context_emission_factors-vehicles.csv"

$ws.Range("A49").Value = "Context Implicit Price Deflators File"
$ws.Range("B49").Value = "String"
$ws.Range("C49").Value = "input_samples/context_implicit_price_deflators.csv"

$ws.Range("A50").Value = "Context Consumer Price Index File"
$ws.Range("B50").Value = "String"
$ws.Range("C50").Value = "input_samples/context_cpi_price_deflators.csv"

# ---------------------------------------------------------------------------
# 3. Sheet view: zoom 90%, freeze pane anchored so row 26 is the first
#    scrollable row, and the active selection sits on C40.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 90
$ws.Range("A26").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("C40").Select()

Write-Host "edit applied"
